$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 29/30: PEPE <-> WrappedeETH swap (name, link) ---
$ws.Range("B29").Value = 'WrappedeETH'
$ws.Range("C29").Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range("B30").Value = 'PEPE'
$ws.Range("C30").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'

# --- Price (D) and Volume(1h) (E) updates ---
$ws.Range("D2").Value = '61.250.91'
$ws.Range("E2").Value = '  -0.99%  '
$ws.Range("D3").Value = '2.390.77'
$ws.Range("E3").Value = '  -4.12%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '548.19'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.31%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '142.28'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -3.58%  '
$ws.Range("E7").Value = '  +0.05%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.539'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -10.86%  '
$ws.Range("D9").Value = '2.389.75'
$ws.Range("E9").Value = '  -4.12%  '
$ws.Range("E10").Value = '  -2.41%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.154'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.19%  '
$ws.Range("E12").Value = '  -2.80%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.348'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -3.63%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '25.37'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -4.10%  '
$ws.Range("D15").Value = '2.823.93'
$ws.Range("E15").Value = '  -3.92%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000166'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -2.33%  '
$ws.Range("D17").Value = '60.938.85'
$ws.Range("E17").Value = '  -1.33%  '
$ws.Range("D18").Value = '2.395.22'
$ws.Range("E18").Value = '  -4.27%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '10.76'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -4.10%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.15'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.92%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '319.11'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.56%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.73'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -4.48%  '
$ws.Range("E23").Value = '  +7.87%  '
$ws.Range("E24").Value = '  +0.06%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '63.89'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.62%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '8.15'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +7.16%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '540.46'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.01%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.00'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.08%  '
$ws.Range("D29").Value = '2.512.77'
$ws.Range("E29").Value = '  -3.76%  '
$ws.Range("D30").Value = '0.0₃0940'
$ws.Range("E30").Value = '  -6.55%  '
$ws.Range("E31").Value = '  -6.40%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '8.12'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -3.78%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.146'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -3.85%  '
$ws.Range("E34").Value = '  -3.50%  '
$ws.Range("E35").Value = '  -1.53%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.999'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.03%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.59'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -6.90%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.73'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -4.59%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.377'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -2.23%  '
$ws.Range("E40").Value = '  +5.93%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '18.14'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -2.66%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '138.96'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -6.42%  '
$ws.Range("E43").Value = '  +0.00%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '40.28'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.34%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.22'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -5.95%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '141.51'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -5.27%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.63'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.61%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '20.29'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -4.45%  '
$ws.Range("E49").Value = '  -3.48%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.579'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -3.68%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0227'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.29%  '
